$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the refreshed cryptos list.
# Numeric-looking Price values must be forced to Text so they stay strings
# (matching the original inlineStr cell type) instead of being auto-converted
# to numbers by Excel.

$ws.Range("D2").Value = "57.060.10"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "2.346.63"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.53"
$ws.Range("E5").Value = "  +2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.76"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").Value = "2.345.81"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("E10").Value = "  +5.71%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("E12").Value = "  +3.02%  "
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.78"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "2.756.69"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "56.868.94"
$ws.Range("E16").Value = "  +3.22%  "
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "2.330.29"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.46"
$ws.Range("E19").Value = "  -2.87%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.55"
$ws.Range("E21").Value = "  +3.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.57"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.84"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("E25").Value = "  +7.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.87"
$ws.Range("E27").Value = "  +4.34%  "
$ws.Range("E28").Value = "  +10.55%  "
$ws.Range("D29").Value = "0.0₃0748"
$ws.Range("E29").Value = "  +4.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.01"
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("E31").Value = "  +5.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.27"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.02"
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("E39").Value = "  +7.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.95"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  +4.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.29"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.21"
$ws.Range("E44").Value = "  +5.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "277.61"
$ws.Range("E45").Value = "  +6.36%  "
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0504"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.84"
$ws.Range("E50").Value = "  +7.26%  "
$ws.Range("E51").Value = "  -0.04%  "

Write-Output "Updated cryptos list"
